# TimeSheet.xlsx update: Alpna time sheet name entry moved to the top of
# the "Alpna" sheet, plus a handful of view/selection/row-height tweaks
# that are left over from browsing through the workbook before saving.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Alpna" sheet: move the filled-in table (header row + 14 entries,
#    previously sitting at A106:F120 underneath a block of blank rows)
#    up to A1:F15, and clear out the old location.
# ---------------------------------------------------------------------
$wsAlpna = $wb.Worksheets.Item("Alpna")
$wsAlpna.Activate()

$src = $wsAlpna.Range("A106:F120")
$src.Copy()
$wsAlpna.Range("A1").PasteSpecial(-4104)
$src.Clear()

$wsAlpna.Range("A1:F15").Select()

# ---------------------------------------------------------------------
# 2. "Deepthi" sheet: selection moved to E1 (scroll position reset).
# ---------------------------------------------------------------------
$wsDeepthi = $wb.Worksheets.Item("Deepthi")
$wsDeepthi.Activate()
$wsDeepthi.Range("E1").Select()

# ---------------------------------------------------------------------
# 3. "Ruchika" sheet: selection moved to G1 (scroll position reset).
# ---------------------------------------------------------------------
$wsRuchika = $wb.Worksheets.Item("Ruchika")
$wsRuchika.Activate()
$wsRuchika.Range("G1").Select()

# ---------------------------------------------------------------------
# 4. "Sravani" sheet: selection moved to F3, and three overly tall rows
#    get shrunk back down.
# ---------------------------------------------------------------------
$wsSravani = $wb.Worksheets.Item("Sravani")
$wsSravani.Activate()
$wsSravani.Range("F3").Select()
$wsSravani.Rows.Item(3).RowHeight = 60
$wsSravani.Rows.Item(16).RowHeight = 75
$wsSravani.Rows.Item(32).RowHeight = 90

# ---------------------------------------------------------------------
# 5. "Sruti" sheet becomes the active tab/sheet when the workbook is
#    saved.
# ---------------------------------------------------------------------
$wsSruti = $wb.Worksheets.Item("Sruti")
$wsSruti.Activate()
